# Apply updated cryptocurrency price/volume data per upstream diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '70.392.83'
$ws.Range("E2").Value = '  -1.90%  '

# Row 3
$ws.Range("D3").Value = '3.607.20'
$ws.Range("E3").Value = '  -0.81%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.05%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '583.34'
$ws.Range("E5").Value = '  -1.27%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '174.26'
$ws.Range("E6").Value = '  -3.25%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.629'
$ws.Range("E7").Value = '  +2.94%  '

# Row 8
$ws.Range("D8").Value = '3.594.87'
$ws.Range("E8").Value = '  -0.86%  '

# Row 9
$ws.Range("E9").Value = '  -0.10%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.193'
$ws.Range("E10").Value = '  -4.57%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.59'
$ws.Range("E11").Value = '  +12.87%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.614'
$ws.Range("E12").Value = '  +1.34%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '48.20'
$ws.Range("E13").Value = '  -3.16%  '

# Row 14
$ws.Range("E14").Value = '  -1.83%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '690.13'
$ws.Range("E15").Value = '  +0.41%  '

# Row 16
$ws.Range("D16").Value = '4.189.43'
$ws.Range("E16").Value = '  -0.72%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '9.00'
$ws.Range("E17").Value = '  +0.17%  '

# Row 18
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '70.396.67'
$ws.Range("E18").Value = '  -2.11%  '

# Row 19
$ws.Range("B19").Value = 'WrappedEther'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D19").Value = '3.606.23'
$ws.Range("E19").Value = '  +0.79%  '

# Row 20
$ws.Range("E20").Value = '  -0.45%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.66'
$ws.Range("E21").Value = '  -3.73%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.38'
$ws.Range("E22").Value = '  -2.07%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.931'
$ws.Range("E23").Value = '  -0.52%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '16.96'
$ws.Range("E24").Value = '  -4.79%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '99.09'
$ws.Range("E25").Value = '  -4.26%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.90'
$ws.Range("E26").Value = '  -3.28%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.75'
$ws.Range("E27").Value = '  -3.44%  '

# Row 28
$ws.Range("E28").Value = '  +0.01%  '

# Row 29
$ws.Range("E29").Value = '  -3.56%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '34.30'
$ws.Range("E30").Value = '  -2.11%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '9.10'
$ws.Range("E31").Value = '  -0.84%  '

# Row 32
$ws.Range("E32").Value = '  -5.44%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.45'
$ws.Range("E33").Value = '  +2.42%  '

# Row 34
$ws.Range("E34").Value = '  -4.90%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.93'
$ws.Range("E35").Value = '  -6.23%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '577.43'
$ws.Range("E36").Value = '  -0.80%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '11.03'
$ws.Range("E37").Value = '  -2.74%  '

# Row 38
$ws.Range("E38").Value = '  -2.30%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '58.28'
$ws.Range("E39").Value = '  -1.91%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.998'
$ws.Range("E40").Value = '  -0.20%  '

# Row 41
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0450'
$ws.Range("E41").Value = '  -3.73%  '

# Row 42
$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.140'
$ws.Range("E42").Value = '  -2.12%  '

# Row 43
$ws.Range("B43").Value = 'TheGraph'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.345'
$ws.Range("E43").Value = '  -1.08%  '

# Row 44
$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").Value = '3.518.23'
$ws.Range("E44").Value = '  -4.31%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '34.06'
$ws.Range("E45").Value = '  -4.86%  '

# Row 46
$ws.Range("D46").Value = '0.0₃0722'
$ws.Range("E46").Value = '  -5.62%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.96'
$ws.Range("E47").Value = '  +5.15%  '

# Row 48
$ws.Range("E48").Value = '  -4.76%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.135'
$ws.Range("E49").Value = '  +1.72%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '135.24'
$ws.Range("E50").Value = '  +2.67%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.85'
$ws.Range("E51").Value = '  -4.82%  '
